$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (resnet50)
$ws.Range("B2").Value = 0.987037037037037
$ws.Range("C2").Value = 0.9865022853476677
$ws.Range("D2").Value = 0.9952380952380953
$ws.Range("E2").Value = 0.9952351171863368

# Row 3 (resnet101)
$ws.Range("B3").Value = 0.9890740740740741
$ws.Range("C3").Value = 0.9887031337958833
$ws.Range("D3").Value = 0.9928571428571429
$ws.Range("E3").Value = 0.9928526757795052

# Row 4 (densenet121)
$ws.Range("B4").Value = 0.9896296296296296
$ws.Range("C4").Value = 0.9893331952174449
$ws.Range("D4").Value = 0.9904761904761905
$ws.Range("E4").Value = 0.9904612778260312

# Row 5 (efficientnet_b0)
$ws.Range("B5").Value = 0.9875925925925926
$ws.Range("C5").Value = 0.9870486448163949

# Row 6 (efficientnet_b3)
$ws.Range("C6").Value = 0.9895197400624811
$ws.Range("D6").Value = 0.9952380952380953
$ws.Range("E6").Value = 0.9952351171863366

# Row 7 (vit_b_16)
$ws.Range("B7").Value = 0.9872222222222222
$ws.Range("C7").Value = 0.9869984774436796
$ws.Range("D7").Value = 0.9809523809523809
$ws.Range("E7").Value = 0.9809002756089702

# Row 8 (swin_t)
$ws.Range("B8").Value = 0.9901851851851852
$ws.Range("C8").Value = 0.9898819612599488
$ws.Range("D8").Value = 0.9976190476190476
$ws.Range("E8").Value = 0.9976175585931684

# Row 9 (convnext_tiny)
$ws.Range("B9").Value = 0.9911111111111112
$ws.Range("C9").Value = 0.9908032995990874
$ws.Range("E9").Value = 0.9952351171863368
